$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-28 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-29 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("755÷8=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "337÷6=56, 1", 2) | Out-Null
$d.Content.Find.Execute("863÷4=215, 3", $true, $false, $false, $false, $false, $true, 1, $false, "972÷3=324, 0", 2) | Out-Null
$d.Content.Find.Execute("327÷7=46, 5", $true, $false, $false, $false, $false, $true, 1, $false, "493÷5=98, 3", 2) | Out-Null
$d.Content.Find.Execute("691÷4=172, 3", $true, $false, $false, $false, $false, $true, 1, $false, "919÷6=153, 1", 2) | Out-Null
$d.Content.Find.Execute("116÷3=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "282÷3=94, 0", 2) | Out-Null
$d.Content.Find.Execute("412÷6=68, 4", $true, $false, $false, $false, $false, $true, 1, $false, "570÷5=114, 0", 2) | Out-Null
$d.Content.Find.Execute("613÷5=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "338÷3=112, 2", 2) | Out-Null
$d.Content.Find.Execute("348÷5=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "672÷9=74, 6", 2) | Out-Null
$d.Content.Find.Execute("186÷5=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "134÷3=44, 2", 2) | Out-Null
$d.Content.Find.Execute("370÷8=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "868÷5=173, 3", 2) | Out-Null
$d.Content.Find.Execute("558÷4=139, 2", $true, $false, $false, $false, $false, $true, 1, $false, "770÷9=85, 5", 2) | Out-Null
$d.Content.Find.Execute("908÷7=129, 5", $true, $false, $false, $false, $false, $true, 1, $false, "121÷9=13, 4", 2) | Out-Null
$d.Content.Find.Execute("785÷7=112, 1", $true, $false, $false, $false, $false, $true, 1, $false, "396÷2=198, 0", 2) | Out-Null
$d.Content.Find.Execute("334÷4=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "992÷8=124, 0", 2) | Out-Null
$d.Content.Find.Execute("644÷6=107, 2", $true, $false, $false, $false, $false, $true, 1, $false, "394÷7=56, 2", 2) | Out-Null
$d.Content.Find.Execute("604÷3=201, 1", $true, $false, $false, $false, $false, $true, 1, $false, "969÷8=121, 1", 2) | Out-Null
$d.Content.Find.Execute("930÷4=232, 2", $true, $false, $false, $false, $false, $true, 1, $false, "882÷8=110, 2", 2) | Out-Null
$d.Content.Find.Execute("604÷2=302, 0", $true, $false, $false, $false, $false, $true, 1, $false, "469÷5=93, 4", 2) | Out-Null
$d.Content.Find.Execute("270÷2=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "135÷4=33, 3", 2) | Out-Null
$d.Content.Find.Execute("725÷5=145, 0", $true, $false, $false, $false, $false, $true, 1, $false, "804÷8=100, 4", 2) | Out-Null
$d.Content.Find.Execute("973÷6=162, 1", $true, $false, $false, $false, $false, $true, 1, $false, "910÷8=113, 6", 2) | Out-Null
$d.Content.Find.Execute("699÷3=233, 0", $true, $false, $false, $false, $false, $true, 1, $false, "327÷9=36, 3", 2) | Out-Null
$d.Content.Find.Execute("856÷3=285, 1", $true, $false, $false, $false, $false, $true, 1, $false, "256÷7=36, 4", 2) | Out-Null
$d.Content.Find.Execute("727÷3=242, 1", $true, $false, $false, $false, $false, $true, 1, $false, "847÷3=282, 1", 2) | Out-Null
$d.Content.Find.Execute("897÷4=224, 1", $true, $false, $false, $false, $false, $true, 1, $false, "468÷8=58, 4", 2) | Out-Null
